$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (and formatting) to the right
$ws.Columns.Item(1).Insert()

# Match the header formatting used by the other header cells (bold, border, centered)
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header for column A
$ws.Cells.Item(1, 1).Value = "Código de pedido"

# Fill in the "Codigo de pedido" values for each data row
# (kept as text, matching how the other code-like columns are stored)
$codigoPedido = @("879436", "887593", "889587", "889583", "889580", "889575", "889562", "889531", "889531", "889364", "889363", "889334", "889217", "889215", "889215", "889212", "889199", "889157", "889150", "889145", "888957", "888818", "888704", "888704", "888704", "888700", "888645", "888641", "888641", "888626", "888626", "888531", "888525", "888525", "888525", "888425", "888424", "889457", "889452", "889452", "888644", "888644", "888642", "888642", "888642", "888634", "888491", "888491", "888490", "888430", "888428", "888428", "888428", "888249", "888244", "888241", "887576")
$colA = $ws.Range("A2:A58")
$colA.NumberFormat = "@"
for ($i = 0; $i -lt $codigoPedido.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $codigoPedido[$i]
}

$wb.Save()
